$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '61.823.06'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '2.413.63'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.06%  '
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").Value = '2.426.66'
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.32%  '
$ws.Range("E13").Value = '  +4.13%  '
$ws.Range("E14").Value = '  +3.76%  '
$ws.Range("E15").Value = '  +6.14%  '
$ws.Range("D16").Value = '2.837.44'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '61.684.76'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").Value = '2.427.50'
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("E23").Value = '  +14.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '622.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.84%  '
$ws.Range("D29").Value = '0.0₃0961'
$ws.Range("E29").Value = '  +5.81%  '
$ws.Range("D30").Value = '2.520.07'
$ws.Range("E31").Value = '  +2.80%  '
$ws.Range("E32").Value = '  +9.42%  '
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '153.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.371'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.36'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.59%  '
$ws.Range("E43").Value = '  +5.30%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").Value = '0.0₆0287'
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.04'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.68%  '
$ws.Range("E50").Value = '  +2.93%  '
$ws.Range("E51").Value = '  +2.96%  '
